$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegTestData")

# Row 2: flip gender from female -> male
$ws.Range("C2").Value = "male"

# Row 3: flip gender from male -> female,
# change alert message to "Please enter first name",
# change status from Pass -> Fail
$ws.Range("C3").Value = "female"
$ws.Range("I3").Value = "Please enter first name"
$ws.Range("J3").Value = "Fail"

# Update the view selection to I3
$ws.Range("I3").Select()
